$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above row 4, shifting Canada/Grecia/Finlandia/Consolidados down.
$ws.Rows("4:6").Insert(-4121)

# Copy formatting from row 3 (Western Sahara) down onto the new rows 4-6.
$ws.Range("A3:L3").Copy()
$ws.Range("A4:L6").PasteSpecial(-4122)

# --- Row 4: Djibouti ---
$ws.Range("A4").Value = "Djibouti"
$ws.Range("B4").Value = 38
$ws.Range("C4").Value = 6656
$ws.Range("E4").Value = 9745
$ws.Range("H4").Value = 8106

# --- Row 5: Qatar ---
$ws.Range("A5").Value = "Qatar"
$ws.Range("B5").Value = 194
$ws.Range("C5").Value = 9352
$ws.Range("E5").Value = 11640
$ws.Range("H5").Value = 11510

# --- Row 6: Uruguay ---
$ws.Range("A6").Value = "Uruguay"
$ws.Range("B6").Value = 734
$ws.Range("C6").Value = 79114
$ws.Range("E6").Value = 99247
$ws.Range("H6").Value = 100147

# Fill formulas across the whole new block at once so Excel stores them as
# shared formulas (matching a fill-down / multi-cell formula entry).
$ws.Range("F4:F6").Formula = "=(E4-C4)/E4"
$ws.Range("I4:I6").Formula = "=(H4-C4)/H4"
$ws.Range("L4:L6").Formula = "=(K4-F4)/K4"

# The row-insert operation above breaks up the F-column shared formula that
# used to span the (now shifted-down) Canada/Grecia/Finlandia rows; restore
# it as a single fill across the shifted range.
$ws.Range("F7:F9").Formula = "=(E7-C7)/E7"

# Match the final view state: scrolled to column B, zoomed to 122%, with
# E8 as the active selected cell.
$ws.Activate()
$excel.ActiveWindow.Zoom = 122
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("E8").Select() | Out-Null

Write-Output "DONE"
